# Horarios 141 update - LP1912 / LP1912-215 / 6203-6173 (15/01/2026)
# Refresh timestamp 08:48:09 -> 08:55:19, new scrape rows appended/merged.

$wb = $excel.ActiveWorkbook

function Set-DataRow($ws, $row) {
    $ws.Cells.Item($row.r, 1).Value = $row.A
    $ws.Cells.Item($row.r, 2).Value = $row.B
    $ws.Cells.Item($row.r, 3).Value = $row.C
    $ws.Cells.Item($row.r, 4).Value = $row.D
    $ws.Cells.Item($row.r, 5).Value = $row.E
}

# --- Sheet "LP1912" ---
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range('A2').Value = 'Última actualización: 08:55:19'
$ws1.Range('A3').Value = 'Total filas: 121'

$ws1Rows = @(
    @{ r=35; A='05:49:10'; B='07:05'; C='23_HERNANDEZ'; D=76; E='LP1912' },
    @{ r=36; A='05:19:24'; B='07:05'; C='15_ABASTO'; D=106; E='LP1912' },
    @{ r=86; A='08:55:19'; B='08:56'; C='16_SANTA ANA'; D=1; E='LP1912' },
    @{ r=87; A='08:55:19'; B='08:56'; C='10_OLMOS'; D=1; E='LP1912' },
    @{ r=88; A='07:19:37'; B='09:01'; C='215A_EL PATO'; D=102; E='LP1912' },
    @{ r=89; A='08:19:33'; B='09:02'; C='23_HERNANDEZ'; D=43; E='LP1912' },
    @{ r=90; A='08:48:09'; B='09:02'; C='215A_EL PATO'; D=14; E='LP1912' },
    @{ r=91; A='07:45:49'; B='09:03'; C='11_ETCHEVERRY'; D=78; E='LP1912' },
    @{ r=92; A='08:48:09'; B='09:04'; C='11_ETCHEVERRY'; D=16; E='LP1912' },
    @{ r=93; A='08:36:20'; B='09:05'; C='23_HERNANDEZ'; D=29; E='LP1912' },
    @{ r=94; A='08:55:19'; B='09:06'; C='23_HERNANDEZ'; D=11; E='LP1912' },
    @{ r=95; A='07:19:37'; B='09:10'; C='16_P MOR-SANTA ANA'; D=111; E='LP1912' },
    @{ r=96; A='08:36:20'; B='09:11'; C='16_SANTA ANA'; D=35; E='LP1912' },
    @{ r=97; A='08:48:09'; B='09:11'; C='16_P MOR-SANTA ANA'; D=23; E='LP1912' },
    @{ r=98; A='08:36:20'; B='09:13'; C='10_OLMOS'; D=37; E='LP1912' },
    @{ r=99; A='08:48:09'; B='09:13'; C='16_SANTA ANA'; D=25; E='LP1912' },
    @{ r=100; A='08:55:19'; B='09:14'; C='16_SANTA ANA'; D=19; E='LP1912' },
    @{ r=101; A='07:19:37'; B='09:16'; C='27_EL RETIRO'; D=117; E='LP1912' },
    @{ r=102; A='07:58:19'; B='09:17'; C='27_EL RETIRO'; D=79; E='LP1912' },
    @{ r=103; A='07:45:49'; B='09:21'; C='26_HERNANDEZ'; D=96; E='LP1912' },
    @{ r=104; A='07:45:49'; B='09:22'; C='17_ROMERO'; D=97; E='LP1912' },
    @{ r=105; A='07:58:19'; B='09:23'; C='17_ROMERO'; D=85; E='LP1912' },
    @{ r=106; A='07:45:49'; B='09:23'; C='11_ETCHEVERRY'; D=98; E='LP1912' },
    @{ r=107; A='08:48:09'; B='09:24'; C='11_ETCHEVERRY'; D=36; E='LP1912' },
    @{ r=108; A='08:19:33'; B='09:25'; C='16_SANTA ANA'; D=66; E='LP1912' },
    @{ r=109; A='07:45:49'; B='09:32'; C='15_ABASTO'; D=107; E='LP1912' },
    @{ r=110; A='07:45:49'; B='09:33'; C='10_OLMOS'; D=108; E='LP1912' },
    @{ r=111; A='08:48:09'; B='09:35'; C='23_HERNANDEZ'; D=47; E='LP1912' },
    @{ r=112; A='07:45:49'; B='09:41'; C='215C_EL PATO'; D=116; E='LP1912' },
    @{ r=113; A='07:58:19'; B='09:42'; C='215C_EL PATO'; D=104; E='LP1912' },
    @{ r=114; A='07:58:19'; B='09:43'; C='14_ABASTO'; D=105; E='LP1912' },
    @{ r=115; A='08:48:09'; B='09:44'; C='14_ABASTO'; D=56; E='LP1912' },
    @{ r=116; A='07:58:19'; B='09:52'; C='15_ABASTO'; D=114; E='LP1912' },
    @{ r=117; A='08:19:33'; B='10:10'; C='16_P MOR-SANTA ANA'; D=111; E='LP1912' },
    @{ r=118; A='08:48:09'; B='10:11'; C='16_P MOR-SANTA ANA'; D=83; E='LP1912' },
    @{ r=119; A='08:19:33'; B='10:12'; C='15_ABASTO'; D=113; E='LP1912' },
    @{ r=120; A='08:36:20'; B='10:21'; C='26_HERNANDEZ'; D=105; E='LP1912' },
    @{ r=121; A='08:36:20'; B='10:22'; C='17_ROMERO'; D=106; E='LP1912' },
    @{ r=122; A='08:36:20'; B='10:26'; C='215A_EL PATO'; D=110; E='LP1912' },
    @{ r=123; A='08:48:09'; B='10:27'; C='215A_EL PATO'; D=99; E='LP1912' },
    @{ r=124; A='08:48:09'; B='10:42'; C='17_ROMERO'; D=114; E='LP1912' },
    @{ r=125; A='08:55:19'; B='10:43'; C='14_ABASTO'; D=108; E='LP1912' },
    @{ r=126; A='08:48:09'; B='10:44'; C='14_ABASTO'; D=116; E='LP1912' }
)
foreach ($row in $ws1Rows) { Set-DataRow $ws1 $row }

# --- Sheet "LP1912-215" (header timestamp only) ---
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range('A2').Value = 'Última actualización: 08:55:19'

# --- Sheet "6203-6173" ---
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range('A2').Value = 'Última actualización: 08:55:19'
$ws3.Range('A3').Value = 'Total filas: 26'

$ws3Rows = @(
    @{ r=25; A='08:55:19'; B='08:55'; C='215A_LA PLATA'; D=0; E='L6173' },
    @{ r=26; A='07:19:37'; B='09:08'; C='215D_LA PLATA'; D=109; E='L6203' },
    @{ r=27; A='07:58:19'; B='09:09'; C='215D_LA PLATA'; D=71; E='L6203' },
    @{ r=28; A='08:36:20'; B='09:10'; C='215D_LA PLATA'; D=34; E='L6203' },
    @{ r=29; A='08:48:09'; B='09:13'; C='215D_LA PLATA'; D=25; E='L6203' },
    @{ r=30; A='08:19:33'; B='10:03'; C='215B_LP-P MOR-40 Y 115'; D=104; E='L6173' },
    @{ r=31; A='08:55:19'; B='10:54'; C='215A_LA PLATA'; D=119; E='L6173' }
)
foreach ($row in $ws3Rows) { Set-DataRow $ws3 $row }

